# Rename field to fieldName / property to PropertyName:
# concretely, add a new "Strength Value" column (E) that mirrors the
# existing "Strength" column (C) values for each creature row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("E1").Value = "Strength Value"

# New column values - same numbers as column C ("Strength") per row
$ws.Range("E2").Value = 5
$ws.Range("E3").Value = 5
$ws.Range("E4").Value = 7
$ws.Range("E5").Value = 9
$ws.Range("E6").Value = 11
$ws.Range("E7").Value = 13
$ws.Range("E8").Value = 15
$ws.Range("E9").Value = 17
$ws.Range("E10").Value = 19
$ws.Range("E11").Value = 21

# Widen the new column similarly to the authored workbook
$ws.Columns.Item(5).ColumnWidth = 12.140625

# Match the authored workbook's view state (zoom + final selection)
$excel.ActiveWindow.Zoom = 200
[void]$ws.Range("G13").Select()
